# Testing GitHub.docx — "modifying the document" edit
#
# The original document is a single paragraph:
#   "Testing GitHub…" followed by the (collapsed) _GoBack bookmark.
#
# The target splits that into four paragraphs:
#   1) "Testing GitHub…"                              (bookmark removed)
#   2) "… modifying the document…" (with a gramStart/gramEnd proofErr
#      pair bracketing the word "modifying", as Word's grammar checker
#      inserts when it flags a sentence start)
#   3) an empty paragraph
#   4) an otherwise-empty paragraph that now owns the _GoBack bookmark

$d = $word.ActiveDocument

# Locate the first (only) paragraph and the bookmark living in it.
$firstPara = $d.Paragraphs.Item(1)
$bm = $d.Bookmarks.Item("_GoBack")

# The bookmark currently sits right before the paragraph mark of
# paragraph 1; remember where the paragraph's text ends (i.e. just
# before that paragraph mark) so we can insert the new content there.
$insertAt = $firstPara.Range.End - 1

# Pull the bookmark out of paragraph 1 — it gets re-added to the new
# trailing paragraph below.
$bm.Delete()

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newXml = "<w:p $w>" +
            "<w:r><w:t xml:space='preserve'>… </w:t></w:r>" +
            "<w:proofErr w:type='gramStart'/>" +
            "<w:r><w:t>modifying</w:t></w:r>" +
            "<w:proofErr w:type='gramEnd'/>" +
            "<w:r><w:t xml:space='preserve'> the document…</w:t></w:r>" +
          "</w:p>" +
          "<w:p $w/>" +
          "<w:p $w/>"

$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.InsertXML($newXml)

# Re-home the _GoBack bookmark onto the new, final (now empty) paragraph.
$lastPara = $d.Paragraphs.Last
$d.Bookmarks.Add("_GoBack", $lastPara.Range)
